# The deck's theme parts (ppt/theme/theme1.xml and ppt/theme/theme2.xml) swap
# their contents: theme1.xml becomes the "Integral" theme (colour scheme +
# name) while theme2.xml becomes the standard "Office Theme" colour scheme.
#
# theme2.xml is the theme actually wired to the slide master / the whole
# deck (ppt/slideMasters/slideMaster1.xml.rels + ppt/_rels/presentation.xml.rels
# both point at theme2.xml), so the net visible effect of the swap is that the
# presentation's live colour scheme changes from the green "Integral" palette
# to the default blue "Office" palette. Apply that palette through the
# PowerPoint theme-colour-scheme object, which writes straight back into the
# shared theme part used by every slide/layout in the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index : Role     : target (Office Theme) colour
#   1   : dark1    : 000000
#   2   : light1   : FFFFFF
#   3   : dark2    : 44546A
#   4   : light2   : E7E6E6
#   5   : accent1  : 5B9BD5
#   6   : accent2  : ED7D31
#   7   : accent3  : A5A5A5
#   8   : accent4  : FFC000
#   9   : accent5  : 4472C4
#  10   : accent6  : 70AD47
#  11   : hyperlink: 0563C1
#  12   : followed : 954F72
# PowerPoint .RGB is an OLE COLORREF packed as R + G*256 + B*65536.
$tcs.Item(1).RGB = 0x000000
$tcs.Item(2).RGB = 0xFFFFFF
$tcs.Item(3).RGB = 0x6A5444
$tcs.Item(4).RGB = 0xE6E6E7
$tcs.Item(5).RGB = 0xD59B5B
$tcs.Item(6).RGB = 0x317DED
$tcs.Item(7).RGB = 0xA5A5A5
$tcs.Item(8).RGB = 0x00C0FF
$tcs.Item(9).RGB = 0xC47244
$tcs.Item(10).RGB = 0x47AD70
$tcs.Item(11).RGB = 0xC16305
$tcs.Item(12).RGB = 0x724F95
